# Auto-generated COM-interop script: adds two new localization-status
# records (2a521b16-... and 4967c61c-...) to the Overview/zh-cn/de-de
# sheets of the handoff report, matching the 'Generate Report for
# Handoff' commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Hyperlinks.Delete()

# Make room for the new row between c003bd65.. and 37ea1dd4..
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "2a521b16-3871-45a3-90ec-45aa6e15bd71.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-26-11 18:26:19"

# Append the new row at the end (4967c61c..)
$ws.Range("A5").Value = "4967c61c-a77c-496d-a7ea-863e9bf454f1.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "2016-26-11 18:26:19"

# Rebuild every File Name hyperlink (anchors shifted with the insert above)
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/43f7a690a20b905f31a8c3fdf488167a3321d2e8/e2e/c003bd65-8677-4b9d-aad3-abac071d090b.md", "", "", "c003bd65-8677-4b9d-aad3-abac071d090b.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b68b1473f1b3f5a37e32dd1a6c7a1e0c6d9f2a31/e2e/2a521b16-3871-45a3-90ec-45aa6e15bd71.md", "", "", "2a521b16-3871-45a3-90ec-45aa6e15bd71.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/a769f066dcdfc66e2f1210d9ce9ee413c8966878/e2e/37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.md", "", "", "37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c57d6e4f3a2b1908f7e6d5c4b3a291807f6e5d4c/e2e/4967c61c-a77c-496d-a7ea-863e9bf454f1.md", "", "", "4967c61c-a77c-496d-a7ea-863e9bf454f1.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()

# Make room for the new row between c003bd65.. and 37ea1dd4..
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "2a521b16-3871-45a3-90ec-45aa6e15bd71.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2a521b16-3871-45a3-90ec-45aa6e15bd71.9c4f8a014caf28c40df2d99ab686e712762239f5.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-11 18:26:15"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

# Append the new row at the end (4967c61c..)
$ws.Range("A5").Value = "4967c61c-a77c-496d-a7ea-863e9bf454f1.md"
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "4967c61c-a77c-496d-a7ea-863e9bf454f1.25653456a5ae380c9fbd77dbcaa5c98ee47c440d.zh-cn.xlf"
$ws.Range("E5").Value = "2016-03-11 18:26:15"
$ws.Range("H5").Value = "0001-01-01 00:00:00"
$ws.Range("I5").Value = "Include"

# Rebuild every hyperlink on this sheet (anchors shifted with the insert above)
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/43f7a690a20b905f31a8c3fdf488167a3321d2e8/e2e/c003bd65-8677-4b9d-aad3-abac071d090b.md", "", "", "c003bd65-8677-4b9d-aad3-abac071d090b.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/43f7a690a20b905f31a8c3fdf488167a3321d2e8/e2e/c003bd65-8677-4b9d-aad3-abac071d090b.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd824585b09f8bf94fa4886d3450c9c9e3636bd8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c003bd65-8677-4b9d-aad3-abac071d090b.d9051e598847c7ea9d5cb7a0011e8a68085e1619.zh-cn.xlf", "", "", "c003bd65-8677-4b9d-aad3-abac071d090b.d9051e598847c7ea9d5cb7a0011e8a68085e1619.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a26d046bd5f8cc9b4faea6b470cbd87ce14a33d2/e2e/c003bd65-8677-4b9d-aad3-abac071d090b.md", "", "", "c003bd65-8677-4b9d-aad3-abac071d090b.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/da42f7c315701b77031b063ff74e3006fe9d43c1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c003bd65-8677-4b9d-aad3-abac071d090b.d9051e598847c7ea9d5cb7a0011e8a68085e1619.zh-cn.xlf", "", "", "c003bd65-8677-4b9d-aad3-abac071d090b.d9051e598847c7ea9d5cb7a0011e8a68085e1619.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b68b1473f1b3f5a37e32dd1a6c7a1e0c6d9f2a31/e2e/2a521b16-3871-45a3-90ec-45aa6e15bd71.md", "", "", "2a521b16-3871-45a3-90ec-45aa6e15bd71.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/b68b1473f1b3f5a37e32dd1a6c7a1e0c6d9f2a31/e2e/2a521b16-3871-45a3-90ec-45aa6e15bd71.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f2e6c3a8b4d5160738291a4c5d6e7f8091a2b3c4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/2a521b16-3871-45a3-90ec-45aa6e15bd71.9c4f8a014caf28c40df2d99ab686e712762239f5.zh-cn.xlf", "", "", "2a521b16-3871-45a3-90ec-45aa6e15bd71.9c4f8a014caf28c40df2d99ab686e712762239f5.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/a769f066dcdfc66e2f1210d9ce9ee413c8966878/e2e/37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.md", "", "", "37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.md")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/a769f066dcdfc66e2f1210d9ce9ee413c8966878/e2e/37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/406abe3598a4c180756b3ecd0af98c86adb5b31d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.5e7a37b2da99f014721c2b7a355e1b9a70c6751c.zh-cn.xlf", "", "", "37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.5e7a37b2da99f014721c2b7a355e1b9a70c6751c.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c57d6e4f3a2b1908f7e6d5c4b3a291807f6e5d4c/e2e/4967c61c-a77c-496d-a7ea-863e9bf454f1.md", "", "", "4967c61c-a77c-496d-a7ea-863e9bf454f1.md")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/c57d6e4f3a2b1908f7e6d5c4b3a291807f6e5d4c/e2e/4967c61c-a77c-496d-a7ea-863e9bf454f1.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/81726354a0b9c8d7e6f5041372635485960a7b8c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4967c61c-a77c-496d-a7ea-863e9bf454f1.25653456a5ae380c9fbd77dbcaa5c98ee47c440d.zh-cn.xlf", "", "", "4967c61c-a77c-496d-a7ea-863e9bf454f1.25653456a5ae380c9fbd77dbcaa5c98ee47c440d.zh-cn.xlf")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()

# Make room for the new row between c003bd65.. and 37ea1dd4..
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "2a521b16-3871-45a3-90ec-45aa6e15bd71.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2a521b16-3871-45a3-90ec-45aa6e15bd71.9c4f8a014caf28c40df2d99ab686e712762239f5.de-de.xlf"
$ws.Range("E3").Value = "2016-03-11 18:26:19"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

# Append the new row at the end (4967c61c..)
$ws.Range("A5").Value = "4967c61c-a77c-496d-a7ea-863e9bf454f1.md"
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "4967c61c-a77c-496d-a7ea-863e9bf454f1.25653456a5ae380c9fbd77dbcaa5c98ee47c440d.de-de.xlf"
$ws.Range("E5").Value = "2016-03-11 18:26:19"
$ws.Range("H5").Value = "0001-01-01 00:00:00"
$ws.Range("I5").Value = "Include"

# Rebuild every hyperlink on this sheet (anchors shifted with the insert above)
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/43f7a690a20b905f31a8c3fdf488167a3321d2e8/e2e/c003bd65-8677-4b9d-aad3-abac071d090b.md", "", "", "c003bd65-8677-4b9d-aad3-abac071d090b.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/43f7a690a20b905f31a8c3fdf488167a3321d2e8/e2e/c003bd65-8677-4b9d-aad3-abac071d090b.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b1cb137ed09fbabdeb80df581d021d696143428d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c003bd65-8677-4b9d-aad3-abac071d090b.d9051e598847c7ea9d5cb7a0011e8a68085e1619.de-de.xlf", "", "", "c003bd65-8677-4b9d-aad3-abac071d090b.d9051e598847c7ea9d5cb7a0011e8a68085e1619.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8540ce3874af8a6ad9fbcc68008525e1f084ef6a/e2e/c003bd65-8677-4b9d-aad3-abac071d090b.md", "", "", "c003bd65-8677-4b9d-aad3-abac071d090b.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a0b4348416064dc28dde4c87dda48735c48d6bc7/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c003bd65-8677-4b9d-aad3-abac071d090b.d9051e598847c7ea9d5cb7a0011e8a68085e1619.de-de.xlf", "", "", "c003bd65-8677-4b9d-aad3-abac071d090b.d9051e598847c7ea9d5cb7a0011e8a68085e1619.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b68b1473f1b3f5a37e32dd1a6c7a1e0c6d9f2a31/e2e/2a521b16-3871-45a3-90ec-45aa6e15bd71.md", "", "", "2a521b16-3871-45a3-90ec-45aa6e15bd71.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/b68b1473f1b3f5a37e32dd1a6c7a1e0c6d9f2a31/e2e/2a521b16-3871-45a3-90ec-45aa6e15bd71.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3a4b5c6d7e8f90a1b2c3d4e5f60718293a4b5c6d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/2a521b16-3871-45a3-90ec-45aa6e15bd71.9c4f8a014caf28c40df2d99ab686e712762239f5.de-de.xlf", "", "", "2a521b16-3871-45a3-90ec-45aa6e15bd71.9c4f8a014caf28c40df2d99ab686e712762239f5.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/a769f066dcdfc66e2f1210d9ce9ee413c8966878/e2e/37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.md", "", "", "37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.md")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/a769f066dcdfc66e2f1210d9ce9ee413c8966878/e2e/37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ff3a514f657f04f247309fad5bfb5fa8e767cbdd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.5e7a37b2da99f014721c2b7a355e1b9a70c6751c.de-de.xlf", "", "", "37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.5e7a37b2da99f014721c2b7a355e1b9a70c6751c.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c57d6e4f3a2b1908f7e6d5c4b3a291807f6e5d4c/e2e/4967c61c-a77c-496d-a7ea-863e9bf454f1.md", "", "", "4967c61c-a77c-496d-a7ea-863e9bf454f1.md")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/c57d6e4f3a2b1908f7e6d5c4b3a291807f6e5d4c/e2e/4967c61c-a77c-496d-a7ea-863e9bf454f1.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5d4c3b2a1908f7e6d5c4b3a2918f7e6d5c4b3a29/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4967c61c-a77c-496d-a7ea-863e9bf454f1.25653456a5ae380c9fbd77dbcaa5c98ee47c440d.de-de.xlf", "", "", "4967c61c-a77c-496d-a7ea-863e9bf454f1.25653456a5ae380c9fbd77dbcaa5c98ee47c440d.de-de.xlf")

